$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values (Date, Volumen, PrecioMin/Max/Prom, Precio $/Kg) for rows 2..23 as they existed BEFORE the edit.
$before = @{
    2  = @{ D = 44469; J = 20; K = 12000; P = 1200 }
    3  = @{ D = 44463; J = 25; K = 12000; P = 1200 }
    4  = @{ D = 44698; J = 35; K = 11000; P = 1100 }
    5  = @{ D = 44656; J = 25; K = 10000; P = 1000 }
    6  = @{ D = 44369; J = 25; K = 8000;  P = 800  }
    7  = @{ D = 44715; J = 30; K = 11000; P = 1100 }
    8  = @{ D = 44525; J = 20; K = 9000;  P = 900  }
    9  = @{ D = 44707; J = 15; K = 12000; P = 1200 }
    10 = @{ D = 44705; J = 30; K = 12000; P = 1200 }
    11 = @{ D = 44473; J = 25; K = 11000; P = 1100 }
    12 = @{ D = 44645; J = 25; K = 10000; P = 1000 }
    13 = @{ D = 44530; J = 30; K = 10000; P = 1000 }
    14 = @{ D = 44425; J = 30; K = 13000; P = 1300 }
    15 = @{ D = 44348; J = 20; K = 10000; P = 1000 }
    16 = @{ D = 44523; J = 30; K = 9000;  P = 900  }
    17 = @{ D = 44649; J = 25; K = 10000; P = 1000 }
    18 = @{ D = 44726; J = 30; K = 14000; P = 1400 }
    19 = @{ D = 44663; J = 30; K = 12000; P = 1200 }
    20 = @{ D = 44722; J = 30; K = 13000; P = 1300 }
    21 = @{ D = 44708; J = 25; K = 11000; P = 1100 }
    22 = @{ D = 44659; J = 25; K = 10000; P = 1000 }
    23 = @{ D = 44526; J = 25; K = 9000;  P = 900  }
}

# Destination row -> source row: the data that now appears in the destination row
# is exactly the data that used to be in the source row (a permutation of the 22 rows).
$mapping = @{
    2  = 18
    3  = 19
    4  = 8
    5  = 20
    6  = 14
    7  = 10
    8  = 5
    9  = 15
    10 = 2
    11 = 16
    12 = 9
    13 = 4
    14 = 3
    15 = 22
    16 = 7
    17 = 12
    18 = 13
    19 = 21
    20 = 17
    21 = 23
    22 = 11
    23 = 6
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $vals = $before[$srcRow]

    $ws.Range("D$destRow").Value = $vals.D
    $ws.Range("J$destRow").Value = $vals.J
    $ws.Range("K$destRow").Value = $vals.K
    $ws.Range("L$destRow").Value = $vals.K
    $ws.Range("M$destRow").Value = $vals.K
    $ws.Range("P$destRow").Value = ($vals.K / 10)
}
